$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product name text on both the input and output sheets
# (the space was dropped from "Late Repayment" -> "LateRepayment").
$newName = "1014-MS-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-LateRepayment"
$ws1.Range("B1").Value = $newName
$ws2.Range("B1").Value = $newName

# Move the selection to B1 on the input sheet (previously A29), then
# select B1 on the output sheet last so it becomes the active sheet/tab.
$ws1.Range("B1").Select()
$ws2.Range("B1").Select()
